# feat: add 2022-Q1 data
#
# The existing "总计" (Total) roll-up sheet becomes the new "2022-Q1"
# per-fund holdings sheet (its old roll-up content is replaced), and a
# fresh "总计" sheet is appended right after it with a new 2022-Q1 row
# prepended to the historical roll-up table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Repurpose the current "总计" sheet as "2022-Q1" (fund holdings).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Numeric-looking columns must stay text (mirrors how every other
# quarter sheet in this workbook stores them), so force Text format
# before writing the values.
$q1.Range("B2:B6").NumberFormat = "@"
$q1.Range("D2:G6").NumberFormat = "@"

$q1.Range("B2").Value = "000179"
$q1.Range("C2").Value = "广发美国房地产指数QDII-人民币"
$q1.Range("D2").Value = "2.37"
$q1.Range("E2").Value = "92.38"
$q1.Range("F2").Value = "4.65"
$q1.Range("G2").Value = "0.1102"
$q1.Range("H2").Value = 3

$q1.Range("B3").Value = "000180"
$q1.Range("C3").Value = "广发美国房地产指数QDII - 美元"
$q1.Range("D3").Value = "2.37"
$q1.Range("E3").Value = "92.38"
$q1.Range("F3").Value = "4.65"
$q1.Range("G3").Value = "0.1102"
$q1.Range("H3").Value = 3

$q1.Range("B4").Value = "160140"
$q1.Range("C4").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)A"
$q1.Range("D4").Value = "1.35"
$q1.Range("E4").Value = "89.10"
$q1.Range("F4").Value = "4.77"
$q1.Range("G4").Value = "0.0644"
$q1.Range("H4").Value = 2

$q1.Range("B5").Value = "070031"
$q1.Range("C5").Value = "嘉实全球房地产(QDII)"
$q1.Range("D5").Value = "0.60"
$q1.Range("E5").Value = "95.08"
$q1.Range("F5").Value = "3.55"
$q1.Range("G5").Value = "0.0213"
$q1.Range("H5").Value = 4

$q1.Range("B6").Value = "160141"
$q1.Range("C6").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)C"
$q1.Range("D6").Value = "0.44"
$q1.Range("E6").Value = "89.10"
$q1.Range("F6").Value = "4.77"
$q1.Range("G6").Value = "0.0210"
$q1.Range("H6").Value = 2

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" roll-up sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.33

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 10
$total.Range("D3").Value = 1.12

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 0.96

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 0.97

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 8
$total.Range("D6").Value = 0.61

$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 7
$total.Range("D7").Value = 0.47

# ---------------------------------------------------------------------
# 3) Match the header / index-column look used throughout the workbook.
#    Cell B1 on "2022-Q1" already carries the bold/centered/bordered
#    header style inherited from the original "总计" sheet, and A2
#    already carries the index-column style - reuse both as the
#    formatting source for every newly created header / index cell.
# ---------------------------------------------------------------------
$q1.Range("B1").Copy()
$q1.Range("C1:H1").PasteSpecial(-4122)
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
